# The commit swaps the content of ppt/theme/theme1.xml (the theme used by
# the presentation's one and only Slide Master - i.e. what the slides
# actually look like) with ppt/theme/theme2.xml (the theme used only by the
# Notes Master). Diffing the "before" theme1.xml/theme2.xml shows their
# <a:fontScheme> and <a:fmtScheme> blocks are already byte-identical - the
# only functional difference is the 12 colours in <a:clrScheme>: theme1.xml
# moves from the custom "Red Violet" palette to the default Office palette
# (the same colours theme2.xml already had).
#
# This COM host has no member that swaps/replaces a whole theme part or the
# Notes Master's theme (Master.Theme / NotesMaster.Theme only expose
# ThemeVariants, which is empty here; *.ApplyTheme(path) is a documented
# no-op stub in this host - it never produces a diff). ThemeColorScheme is
# the one theme-editing surface that is actually wired up and does persist
# into ppt/theme/theme1.xml, so we use it to rewrite the 12 scheme colours
# on the (only) Slide Master's theme to reproduce the visible effect of the
# swap.

function Convert-HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" colours, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$targetHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = Convert-HexToRgbInt $targetHex[$i - 1]
}
